# Weekly fruit/vegetable price update.
# A new observation (row) is inserted at row 72 ("Ají" / Agrícola del Norte
# S.A. de Arica), pushing the existing rows 72-103 down to 73-104 and
# extending the used range from A1:R103 to A1:R104. The newly inserted
# row 72 is then populated with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 72, shifting rows 72..103 down to 73..104.
$ws.Rows.Item(72).Insert()

# Fill the newly inserted row 72 with this week's record.
$ws.Cells.Item(72, 1).Value = 1
$ws.Cells.Item(72, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(72, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(72, 4).Value = 44830
$ws.Cells.Item(72, 5).Value = 15
$ws.Cells.Item(72, 6).Value = 100112021
$ws.Cells.Item(72, 7).Value = "Ají"
$ws.Cells.Item(72, 8).Value = "Inferno"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 160
$ws.Cells.Item(72, 11).Value = 16000
$ws.Cells.Item(72, 12).Value = 17000
$ws.Cells.Item(72, 13).Value = 16500
$ws.Cells.Item(72, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 1100
$ws.Cells.Item(72, 17).Value = 15
$ws.Cells.Item(72, 18).Value = "Hortaliza"
